$wb = $excel.ActiveWorkbook

# --- Sheet "Ingreso": append 9 new aporte rows (2023-06-24 / serial 45101) ---
$ws1 = $wb.Worksheets.Item("Ingreso")

$ingresoRows = @(
    @{ Row = 455; Nombre = "Invitados"; Monto = 200 },
    @{ Row = 456; Nombre = "Michy";     Monto = 100 },
    @{ Row = 457; Nombre = "Randy";     Monto = 100 },
    @{ Row = 458; Nombre = "Anuel";     Monto = 100 },
    @{ Row = 459; Nombre = "Punto";     Monto = 250 },
    @{ Row = 460; Nombre = "Johan";     Monto = 100 },
    @{ Row = 461; Nombre = "Jeicol";    Monto = 100 },
    @{ Row = 462; Nombre = "Omaury";    Monto = 100 },
    @{ Row = 463; Nombre = "Javier";    Monto = 600 }
)

foreach ($r in $ingresoRows) {
    $row = $r.Row
    $ws1.Cells.Item($row, 1).Value = 45101
    $ws1.Cells.Item($row, 2).Value = $r.Nombre
    $ws1.Cells.Item($row, 3).Value = $r.Monto
    $ws1.Cells.Item($row, 4).Value = "Aporte"
}
# The very first new row inherits the column's "Aporte monto" number style,
# matching the prior row; the rest were typed without that style carried over.
for ($row = 456; $row -le 463; $row++) {
    $ws1.Cells.Item($row, 3).Style = "Normal"
}

# --- Sheet "Gastos": append 1 new gasto row ---
$ws2 = $wb.Worksheets.Item("Gastos")
$ws2.Cells.Item(48, 1).Value = 45101
$ws2.Cells.Item(48, 2).Value = "Arbitro, agua y hielo"
$ws2.Cells.Item(48, 3).Value = 940

# --- Update view/selection state to match where the author ended up ---
[void]$ws2.Activate()
[void]$ws2.Range("A48").Select()
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 29
$win2.ScrollColumn = 1

[void]$ws1.Activate()
[void]$ws1.Range("C460").Select()
$win1 = $excel.ActiveWindow
$win1.ScrollRow = 442
$win1.ScrollColumn = 1
